# stocks.xlsx - "updates in the readme version 2.31"
# Replace the ticker list with a single ADBE row, blank out the rest of
# the old symbols, and mark the (now empty) header/quote row as
# explicit text (quote-prefixed, black font) per the refreshed template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: DIS -> ADBE
$ws.Range("A2").Value = "ADBE"

# A3: used to hold "SONY"; the template now keeps this row but forces it
# to an empty, quote-prefixed text cell (black font instead of the theme
# color) and gives the row its new height.
$ws.Range("A3").Value = "'"
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").RowHeight = 19.5

# A4:A7 (NIO, META, NVDA, TSLA) are cleared out entirely.
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = ""
